# Submit work for checkout form - update availability successfully
#
# Changes:
#  - Coding sheet: C3 task owner reassigned from Jeffrey to Tin
#  - Coding sheet: new "x" availability markers added in column E for rows 2-7
#    (boxed with a thin left/right border for rows 2-5, and a thin left
#    border only for rows 6-7)
#  - Coding sheet: row 5 gets an additional note in column F ("without
#    authors"), boxed the same way as the rest of column E on that row
#  - Coding sheet: new narrow column E added, selection moved to E1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coding")

# --- Reassign the "Login form" task owner (Jeffrey -> Tin) ---
$ws.Range("C3").Value = "Tin"

# --- New "x" marks in column E for rows 2-7, "without authors" note in F5 ---
$ws.Range("E2").Value = "x"
$ws.Range("E3").Value = "x"
$ws.Range("E4").Value = "x"
$ws.Range("E5").Value = "x"
$ws.Range("F5").Value = "without authors"
$ws.Range("E6").Value = "x"
$ws.Range("E7").Value = "x"

# --- Borders: box rows 2-5 (E+F) with a thin left/right border ---
$ws.Range("E2:E5").Borders.Item(7).LineStyle = 1
$ws.Range("E2:E5").Borders.Item(7).Weight = 2
$ws.Range("E2:E5").Borders.Item(10).LineStyle = 1
$ws.Range("E2:E5").Borders.Item(10).Weight = 2

$ws.Range("F5").Borders.Item(7).LineStyle = 1
$ws.Range("F5").Borders.Item(7).Weight = 2
$ws.Range("F5").Borders.Item(10).LineStyle = 1
$ws.Range("F5").Borders.Item(10).Weight = 2

# --- Borders: rows 6-7 only get a thin left border ---
$ws.Range("E6:E7").Borders.Item(7).LineStyle = 1
$ws.Range("E6:E7").Borders.Item(7).Weight = 2

# --- Narrow column E width, matching the new marker column ---
$ws.Columns.Item(5).ColumnWidth = 1.53125

# --- Move the active selection to E1, like the source workbook ---
$ws.Range("E1").Select()
